$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Year" -> "Dates"
$ws.Range("A1").Value = "Dates"

# New date labels (stored as literal text, not Excel dates) and new GDP values
$dates = @(
    "2023-10-01",
    "2024-01-01",
    "2024-04-01",
    "2024-07-01",
    "2024-10-01",
    "2025-01-01",
    "2025-04-01",
    "2025-07-01",
    "2025-10-01",
    "2026-01-01",
    "2026-04-01",
    "2026-07-01"
)
$gdp = @(
    -0.02,
    -0.06,
    0.23,
    -0.07000000000000001,
    0.01,
    0.08,
    0.14,
    0.31,
    0.14,
    -0.15,
    -0.25,
    0.08
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2

    $cellA = $ws.Cells.Item($row, 1)
    # Leading apostrophe forces text interpretation (prevents the date-like
    # string from being re-parsed into a date serial number).
    $cellA.Value = "'" + $dates[$i]
    # Drop the style back to the workbook default so the cell carries no
    # explicit format (matches the original "Year" column's unstyled cells).
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $gdp[$i]
}
